# Week 15 simulations update
# Updates rushing and receiving stat cells for the Bengals "Players Data" workbook.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$wsRushing = $wb.Worksheets.Item("Rushing")

# Row 2 - J.Burrow
$wsRushing.Range("D2").Value = 9

# Row 4 - J.Mixon
$wsRushing.Range("C4").Value = 151
$wsRushing.Range("D4").Value = 78
$wsRushing.Range("F4").Value = 29

# Row 5 - S.Perine
$wsRushing.Range("C5").Value = 21
$wsRushing.Range("D5").Value = 19
$wsRushing.Range("E5").Value = 10

# Row 8 - T.Boyd
$wsRushing.Range("C8").Value = 2

# Row 9 - J.Chase
$wsRushing.Range("C9").Value = 5

# --- Receiving sheet ---
$wsReceiving = $wb.Worksheets.Item("Receiving")

# Row 2 - J.Mixon
$wsReceiving.Range("C2").Value = 31
$wsReceiving.Range("D2").Value = 26

# Row 3 - S.Perine
$wsReceiving.Range("C3").Value = 25
$wsReceiving.Range("D3").Value = 21
$wsReceiving.Range("G3").Value = 5
$wsReceiving.Range("H3").Value = 5

# Row 5 - T.Boyd
$wsReceiving.Range("C5").Value = 62
$wsReceiving.Range("D5").Value = 47
$wsReceiving.Range("E5").Value = 16
$wsReceiving.Range("F5").Value = 9

# Row 6 - T.Higgins
$wsReceiving.Range("C6").Value = 74
$wsReceiving.Range("D6").Value = 52
$wsReceiving.Range("E6").Value = 27
$wsReceiving.Range("F6").Value = 16

# Row 7 - J.Chase
$wsReceiving.Range("C7").Value = 68
$wsReceiving.Range("D7").Value = 46
$wsReceiving.Range("E7").Value = 31
$wsReceiving.Range("F7").Value = 15
$wsReceiving.Range("G7").Value = 10
$wsReceiving.Range("H7").Value = 6

# Row 11 - C.Uzomah
$wsReceiving.Range("C11").Value = 38
$wsReceiving.Range("D11").Value = 30
$wsReceiving.Range("E11").Value = 8
$wsReceiving.Range("F11").Value = 7

# Row 12 - D.Sample
$wsReceiving.Range("C12").Value = 10
$wsReceiving.Range("D12").Value = 7
